# Logged 2021 divisional round, simulated season from conference round.
# Appends the new game's per-play yardage logs / special-teams logs to the
# running season strings, and bumps the season-total numeric stat cells by
# the new game's contribution.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# YDS sheet: per-play yardage logs (space separated numbers), append the
# divisional-round game's plays to each existing log.
# ---------------------------------------------------------------------
$ws = $wb.Sheets.Item("YDS")

$ws.Range("B2").Value = $ws.Range("B2").Value() + " 4 -1 1 5 6 14 3 4 4 3 2 -1 2 3 2 5 4 3 2 2"
$ws.Range("B3").Value = $ws.Range("B3").Value() + " 14 19 10 11 7 9 1 3 8 9 5 75 6 6 25 1 5 6 2 4"
$ws.Range("C2").Value = $ws.Range("C2").Value() + " 1 2 5 -1 2 2 7 4 2 6 4 5 4 5 1 5 9 5 8 13 5 -2 4 3 0 9 1 1"
$ws.Range("C3").Value = $ws.Range("C3").Value() + " 15 16 12 18 3 12 -1 24 6 12 14"

# ---------------------------------------------------------------------
# OFF sheet: season offensive totals, row 2 = Home, row 3 = Road.
# ---------------------------------------------------------------------
$ws = $wb.Sheets.Item("OFF")

$ws.Range("C2").Value = 220
$ws.Range("F2").Value = 56
$ws.Range("G2").Value = 73
$ws.Range("J2").Value = 26
$ws.Range("L2").Value = 321
$ws.Range("M2").Value = 224
$ws.Range("Q2").Value = 580

$ws.Range("B3").Value = 13
$ws.Range("C3").Value = 182
$ws.Range("E3").Value = 25
$ws.Range("F3").Value = 118
$ws.Range("H3").Value = 29
$ws.Range("I3").Value = 62
$ws.Range("J3").Value = 62
$ws.Range("N3").Value = 19

# ---------------------------------------------------------------------
# DEF sheet: season defensive totals, row 2 = Home, row 3 = Road.
# ---------------------------------------------------------------------
$ws = $wb.Sheets.Item("DEF")

$ws.Range("C2").Value = 188
$ws.Range("E2").Value = 4
$ws.Range("F2").Value = 56
$ws.Range("G2").Value = 52
$ws.Range("I2").Value = 6
$ws.Range("J2").Value = 36
$ws.Range("L2").Value = 318
$ws.Range("M2").Value = 192
$ws.Range("O2").Value = 31
$ws.Range("Q2").Value = 539

$ws.Range("B3").Value = 12
$ws.Range("C3").Value = 185
$ws.Range("E3").Value = 33
$ws.Range("F3").Value = 114
$ws.Range("G3").Value = 46
$ws.Range("H3").Value = 20
$ws.Range("I3").Value = 57
$ws.Range("J3").Value = 62
$ws.Range("N3").Value = 25

# ---------------------------------------------------------------------
# ST sheet: special teams - kickoff (B col) / punt (D col) logs + totals.
# ---------------------------------------------------------------------
$ws = $wb.Sheets.Item("ST")

$ws.Range("B2").Value = 90
$ws.Range("D2").Value = 55
$ws.Range("F2").Value = 716
$ws.Range("G2").Value = 699
$ws.Range("J2").Value = 260
$ws.Range("K2").Value = 236

$ws.Range("D3").Value = $ws.Range("D3").Value() + " 42 59 36 41 57"
$ws.Range("B4").Value = $ws.Range("B4").Value() + " 57 60 51"
$ws.Range("D4").Value = $ws.Range("D4").Value() + " 0 6 0 7 9"
$ws.Range("B5").Value = $ws.Range("B5").Value() + " 32 45 14"
$ws.Range("D5").Value = $ws.Range("D5").Value() + " 0 0 6 5 0"
$ws.Range("B6").Value = $ws.Range("B6").Value() + " 15 26 25"

# ---------------------------------------------------------------------
# TURNS sheet: turnovers.
# ---------------------------------------------------------------------
$ws = $wb.Sheets.Item("TURNS")

$ws.Range("C2").Value = 13
$ws.Range("D2").Value = 8

# ---------------------------------------------------------------------
# PEN sheet: penalties.
# ---------------------------------------------------------------------
$ws = $wb.Sheets.Item("PEN")

$ws.Range("B2").Value = 12
